$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.245.05"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.356.44"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.92"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.96"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.506"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.21"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.73"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.721.26"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.365.13"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.223.82"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.24"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  +4.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.23"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.37"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.56"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.14"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.48"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.03"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0726"
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.24"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.41"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.63"
$ws.Range("E39").Value = "  +13.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.43"
$ws.Range("E42").Value = "  -34.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.945.50"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.48"
$ws.Range("E46").Value = "  -9.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.584.77"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.97"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.17"
$ws.Range("E51").Value = "  +1.18%  "
